$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows (for team members L. James / L. Messi / K. Bryant)
# right after the existing "A" team roster (before row 4), pushing the
# schedule block and the "Modele" block down by 3 rows.
$ws.Rows("4:6").Insert()

# --- Row 1: team header letters for the 3 new teams (B, C, D), matching
# the existing bold/black/centered/wrap style used by B1 ("A"). Copy B1's
# format onto the new cells instead of rebuilding it property-by-property
# so no redundant style entries are left behind in the stylesheet.
$ws.Range("C1").Value = "B"
$ws.Range("D1").Value = "C"
$ws.Range("E1").Value = "D"
$ws.Range("B1").Copy()
$ws.Range("C1:E1").PasteSpecial(-4122)

# --- Team rosters (6 players per team, 4 teams: A, B, C, D)
$ws.Range("C2").Value = "T. Woods"
$ws.Range("C3").Value = "R. Federer"
$ws.Range("C4").Value = "P. Mickelson"
$ws.Range("C5").Value = "R. Nadal"
$ws.Range("C6").Value = "M. Ryan"

$ws.Range("D2").Value = "M. Pacquiao"
$ws.Range("D3").Value = "Z. Ibrahimović"
$ws.Range("D4").Value = "D. Rose"
$ws.Range("D5").Value = "G. Bale"
$ws.Range("D6").Value = "R. Falcao"

$ws.Range("E2").Value = "M. Özil"
$ws.Range("E3").Value = "N. Djokovic"
$ws.Range("E4").Value = "M. Stafford"
$ws.Range("E5").Value = "L. Hamilton"
$ws.Range("E6").Value = "K. Durant"

$ws.Range("B4").Value = "L. James"
$ws.Range("B5").Value = "L. Messi"
$ws.Range("B6").Value = "K. Bryant"

# --- Schedule section header ("Semaine 14" replaces "Semaine 1")
$ws.Range("A9").Value = "Semaine 14"

# --- Day labels with the new dates (week of 2022-04-04)
$ws.Range("B10").Value = "Lundi" + [char]10 + "2022-04-04"
$ws.Range("B11").Value = "Mardi" + [char]10 + "2022-04-05"
$ws.Range("B12").Value = "Mercredi" + [char]10 + "2022-04-06"
$ws.Range("B13").Value = "Jeudi" + [char]10 + "2022-04-07"
$ws.Range("B14").Value = "Vendredi" + [char]10 + "2022-04-08"
$ws.Range("B15").Value = "Samedi" + [char]10 + "2022-04-09"
$ws.Range("B16").Value = "dimanche" + [char]10 + "2022-04-10"
$ws.Rows("10:16").AutoFit()

# --- Shift assignments for the week (Q1/Q2/Q3 columns)
$ws.Range("C10").Value = "A B"
$ws.Range("D10").Value = "C D"
$ws.Range("E10").Value = "A"

$ws.Range("C11").Value = "B C"
$ws.Range("D11").Value = "D"
$ws.Range("E11").Value = "A B"

$ws.Range("C12").Value = "C D"
$ws.Range("D12").Value = "A"
$ws.Range("E12").Value = "B C"

$ws.Range("C13").Value = "D"
$ws.Range("D13").Value = "A B"
$ws.Range("E13").Value = "C D"

# --- New "Modele" summary block
$ws.Range("A19").Value = "Modele : repartition concentre h-pers = 1050"

$notes = " date de réf. :2022-04-08 12:12:00 " + [char]10 + " sem : 14" + [char]10 + " Modele : repartition concentre h-pers = 1050" + [char]10 + " Calcul présences totales d'équipes: 30" + [char]10 + " Calcul présences individuelles: 150.0" + [char]10 + " Créneaux par jour: 3" + [char]10 + " Equipes par créneau: 2" + [char]10 + " Nombre d'équipes: 4" + [char]10 + " Empl. par éq.: 5" + [char]10 + " Durée quart.: 7.0" + [char]10
$ws.Range("A20").Value = $notes
$ws.Rows("19:20").AutoFit()
